$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark from its original spot (right after
#    "... pregled pojedinačnog studenta je uspešno izvršen" and before the
#    closing period). Plain bookmark deletion does not coalesce the runs
#    that sit on either side of it.
# ---------------------------------------------------------------------------
$oldMark = $d.Bookmarks.Item("_GoBack")
$oldMark.Delete()

# ---------------------------------------------------------------------------
# 2) Locate the run that currently reads
#      " na sistemu. Sve što treba da odradi je da na "
#    and rewrite it as
#      " na sistemu bez opcije da menja detalje date role. Sve što treba da
#        odradi je da na "
#    while keeping it from merging with its left ("role za nekog
#    registrovanog studenta") or right ("stranici sa ...") neighbour runs.
#    We do that by temporarily bookmarking both boundaries of the run before
#    touching its text, then removing those helper bookmarks again.
# ---------------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute(" na sistemu. Sve što treba da odradi je da na ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$leftPin = $d.Range($target.Start, $target.Start)
$d.Bookmarks.Add("ZZZLEFTPIN", $leftPin) | Out-Null
$rightPin = $d.Range($target.End, $target.End)
$d.Bookmarks.Add("ZZZRIGHTPIN", $rightPin) | Out-Null

$target.Text = " na sistemu bez opcije da menja detalje date role. Sve što treba da odradi je da na "

$d.Bookmarks.Item("ZZZLEFTPIN").Delete()
$d.Bookmarks.Item("ZZZRIGHTPIN").Delete()

# ---------------------------------------------------------------------------
# 3) Split that single run into three runs (identical formatting) with an
#    empty "_GoBack" bookmark sitting between run #2 and run #3:
#      " na sistemu" | " bez opcije da menja detalje date role" | _GoBack |
#      ". Sve što treba da odradi je da na "
#    Wrapping the middle chunk in a (temporary) bookmark forces the engine
#    to keep it as a separate run on save; deleting that temporary bookmark
#    afterwards leaves the split in place, so we then drop the real empty
#    "_GoBack" bookmark exactly at the boundary it vacates.
# ---------------------------------------------------------------------------
$middle = $d.Content
$middle.Find.Execute(" bez opcije da menja detalje date role", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $middle.End

$d.Bookmarks.Add("ZZZWRAP", $middle) | Out-Null
$d.Bookmarks.Item("ZZZWRAP").Delete()

$goBackRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
